$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 354.3
$ws.Range("I135").Value = 257.16666
$ws.Range("K135").Value = 2314.49994
$ws.Range("M135").Value = 220.5000600000003
$ws.Range("H138").Value = 506839.12
$ws.Range("I138").Value = 685.2778
$ws.Range("K138").Value = 2055.8334
$ws.Range("M138").Value = 3084.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7851
$ws.Range("I2").Value = 787.5
$ws.Range("J2").Value = 25509.75
$ws.Range("K2").Value = 787.5
$ws.Range("L2").Value = 25509.75
$ws.Range("M2").Value = -674.5
$ws.Range("N2").Value = -25735.75
$ws.Range("H32").Value = 4872.645
$ws.Range("I32").Value = 5716.92
$ws.Range("J32").Value = 1354.8334
$ws.Range("K32").Value = 5716.92
$ws.Range("L32").Value = 1354.8334
$ws.Range("M32").Value = -5429.92
$ws.Range("N32").Value = -1928.8334
$ws.Range("H61").Value = 1611
$ws.Range("I61").Value = 1152.625
$ws.Range("J61").Value = 2833.3333
$ws.Range("K61").Value = 1152.625
$ws.Range("L61").Value = 2833.3333
$ws.Range("M61").Value = -940.625
$ws.Range("N61").Value = -3257.3333
$ws.Range("H74").Value = 1250.0769
$ws.Range("I74").Value = 851.625
$ws.Range("K74").Value = 851.625
$ws.Range("M74").Value = 22.375
$ws.Range("H77").Value = 1250.0769
$ws.Range("I77").Value = 851.625
$ws.Range("K77").Value = 4258.125
$ws.Range("M77").Value = 109.875
$ws.Range("H116").Value = 7851
$ws.Range("I116").Value = 787.5
$ws.Range("J116").Value = 25509.75
$ws.Range("K116").Value = 787.5
$ws.Range("L116").Value = 25509.75
$ws.Range("M116").Value = 1506.5
$ws.Range("N116").Value = -30097.75
$ws.Range("H132").Value = 1878.721
$ws.Range("I132").Value = 1504.4546
$ws.Range("K132").Value = 4513.3638
$ws.Range("M132").Value = -1983.3638
$ws.Range("H136").Value = 1611
$ws.Range("I136").Value = 1152.625
$ws.Range("J136").Value = 2833.3333
$ws.Range("K136").Value = 3457.875
$ws.Range("L136").Value = 8499.999899999999
$ws.Range("M136").Value = -907.875
$ws.Range("N136").Value = -13599.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7851
$ws.Range("I3").Value = 787.5
$ws.Range("J3").Value = 25509.75
$ws.Range("K3").Value = 787.5
$ws.Range("L3").Value = 25509.75
$ws.Range("M3").Value = -673.5
$ws.Range("N3").Value = -25737.75
$ws.Range("H68").Value = 16500
$ws.Range("J68").Value = 16500
$ws.Range("L68").Value = 16500
$ws.Range("N68").Value = -18122
$ws.Range("H71").Value = 16500
$ws.Range("J71").Value = 16500
$ws.Range("L71").Value = 49500
$ws.Range("N71").Value = -57612

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 100001320
$ws.Range("I16").Value = 142858400
$ws.Range("J16").Value = 1470.6666
$ws.Range("K16").Value = 142858400
$ws.Range("L16").Value = 1470.6666
$ws.Range("M16").Value = -142858113
$ws.Range("N16").Value = -2044.6666
$ws.Range("H31").Value = 2249.8125
$ws.Range("I31").Value = 1039.7
$ws.Range("K31").Value = 1039.7
$ws.Range("M31").Value = -744.7
$ws.Range("H34").Value = 2249.8125
$ws.Range("I34").Value = 1039.7
$ws.Range("K34").Value = 1039.7
$ws.Range("M34").Value = -837.7
$ws.Range("H58").Value = 1126.3529
$ws.Range("I58").Value = 959.2
$ws.Range("J58").Value = 1365.1428
$ws.Range("K58").Value = 959.2
$ws.Range("L58").Value = 1365.1428
$ws.Range("M58").Value = -756.2
$ws.Range("N58").Value = -1771.1428
$ws.Range("H109").Value = 23612.625
$ws.Range("J109").Value = 23612.625
$ws.Range("L109").Value = 23612.625
$ws.Range("N109").Value = -25692.625
$ws.Range("H113").Value = 100001320
$ws.Range("I113").Value = 142858400
$ws.Range("J113").Value = 1470.6666
$ws.Range("K113").Value = 142858400
$ws.Range("L113").Value = 1470.6666
$ws.Range("M113").Value = -142856230
$ws.Range("N113").Value = -5810.6666
$ws.Range("H132").Value = 6102.087
$ws.Range("I132").Value = 6861.7646
$ws.Range("J132").Value = 3949.6667
$ws.Range("K132").Value = 20585.2938
$ws.Range("L132").Value = 11849.0001
$ws.Range("M132").Value = -18055.2938
$ws.Range("N132").Value = -16909.0001
$ws.Range("H134").Value = 25642752
$ws.Range("I134").Value = 27779480
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 83338440
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -83335905
$ws.Range("N134").Value = -11070
$ws.Range("H136").Value = 1126.3529
$ws.Range("I136").Value = 959.2
$ws.Range("J136").Value = 1365.1428
$ws.Range("K136").Value = 2877.6
$ws.Range("L136").Value = 4095.4284
$ws.Range("M136").Value = -327.6000000000004
$ws.Range("N136").Value = -9195.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1557
$ws.Range("J34").Value = 1935.7142
$ws.Range("L34").Value = 5807.142599999999
$ws.Range("N34").Value = -5975.142599999999
$ws.Range("H131").Value = 21278006
$ws.Range("J131").Value = 1566.2195
$ws.Range("L131").Value = 4698.6585
$ws.Range("N131").Value = -14778.6585

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1933.3846
$ws.Range("J113").Value = 2464.4
$ws.Range("L113").Value = 2464.4
$ws.Range("N113").Value = -6804.4
$ws.Range("H122").Value = 1302.3334
$ws.Range("I122").Value = 1302.3334
$ws.Range("K122").Value = 3907.0002
$ws.Range("M122").Value = -1457.0002
$ws.Range("H132").Value = 2131.2307
$ws.Range("I132").Value = 1777.2727
$ws.Range("K132").Value = 5331.8181
$ws.Range("M132").Value = -2801.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 286956.66
$ws.Range("I2").Value = 400696
$ws.Range("J2").Value = 205714.28
$ws.Range("K2").Value = 400696
$ws.Range("L2").Value = 205714.28
$ws.Range("M2").Value = -400584
$ws.Range("N2").Value = -205938.28
$ws.Range("H87").Value = 18500
$ws.Range("I87").Value = 15000
$ws.Range("J87").Value = 22000
$ws.Range("K87").Value = 15000
$ws.Range("L87").Value = 22000
$ws.Range("M87").Value = -13877
$ws.Range("N87").Value = -24246
$ws.Range("H88").Value = 10500
$ws.Range("I88").Value = 10500
$ws.Range("K88").Value = 10500
$ws.Range("M88").Value = -10072
$ws.Range("H90").Value = 18500
$ws.Range("I90").Value = 15000
$ws.Range("J90").Value = 22000
$ws.Range("K90").Value = 45000
$ws.Range("L90").Value = 66000
$ws.Range("M90").Value = -39384
$ws.Range("N90").Value = -77232
$ws.Range("H91").Value = 10500
$ws.Range("I91").Value = 10500
$ws.Range("K91").Value = 10500
$ws.Range("M91").Value = -9018
$ws.Range("H132").Value = 69112.266
$ws.Range("I132").Value = 2309.5
$ws.Range("J132").Value = 145458.28
$ws.Range("K132").Value = 6928.5
$ws.Range("L132").Value = 436374.84
$ws.Range("M132").Value = -4398.5
$ws.Range("N132").Value = -441434.84
$ws.Range("H136").Value = 10067.667
$ws.Range("I136").Value = 12501.333
$ws.Range("K136").Value = 37503.999
$ws.Range("M136").Value = -34953.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10121.667
$ws.Range("J41").Value = 10121.667
$ws.Range("L41").Value = 10121.667
$ws.Range("N41").Value = -10901.667
$ws.Range("H45").Value = 8729.200000000001
$ws.Range("J45").Value = 8729.200000000001
$ws.Range("L45").Value = 8729.200000000001
$ws.Range("N45").Value = -9711.200000000001
$ws.Range("H74").Value = 8845.200000000001
$ws.Range("J74").Value = 8845.200000000001
$ws.Range("L74").Value = 8845.200000000001
$ws.Range("N74").Value = -10717.2
$ws.Range("H77").Value = 8845.200000000001
$ws.Range("J77").Value = 8845.200000000001
$ws.Range("L77").Value = 26535.6
$ws.Range("N77").Value = -35895.60000000001
$ws.Range("H107").Value = 512.44446
$ws.Range("I107").Value = 599
$ws.Range("J107").Value = 443.2
$ws.Range("K107").Value = 1797
$ws.Range("L107").Value = 1329.6
$ws.Range("M107").Value = 123
$ws.Range("N107").Value = -5169.6
$ws.Range("H115").Value = 47332.5
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 47332.5
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 47332.5
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -50466.5
$ws.Range("H132").Value = 5074.5713
$ws.Range("I132").Value = 4746
$ws.Range("J132").Value = 5666
$ws.Range("K132").Value = 14238
$ws.Range("L132").Value = 16998
$ws.Range("M132").Value = -11708
$ws.Range("N132").Value = -22058
$ws.Range("H135").Value = 43223
$ws.Range("J135").Value = 43223
$ws.Range("L135").Value = 43223
$ws.Range("N135").Value = -53363
$ws.Range("H136").Value = 858.8
$ws.Range("I136").Value = 512.5714
$ws.Range("J136").Value = 1666.6666
$ws.Range("K136").Value = 1537.7142
$ws.Range("L136").Value = 4999.9998
$ws.Range("M136").Value = 1012.2858
$ws.Range("N136").Value = -10099.9998
